# Add a new "Gujarat DC" location row (row 14) to the Setup sheet, mirroring
# the structure of the existing row 13 ("Chennai DC"), and update the
# selection to reflect where the user ended up after data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")

# Copy formatting (styles) cell-by-cell from row 13 into row 14 so the new
# row reuses the existing style indices instead of minting new ones.
$cols = @("A", "B", "D", "E", "F", "G", "K", "N", "O", "P")
foreach ($col in $cols) {
    $ws.Range($col + "13").Copy()
    $ws.Range($col + "14").PasteSpecial(-4122)
}

# Populate the new row's values.
$ws.Range("A14").Value = "IBM APAC"
$ws.Range("B14").Value = "Classification"
$ws.Range("D14").Value = "G20-LocationService"
$ws.Range("E14").Value = "G20-LocationService-L2"
$ws.Range("F14").Value = "G20-LocationService-L3"
$ws.Range("G14").Value = "G20-LS- Gujarat DC"
$ws.Range("K14").Value = "Gujarat"
$ws.Range("N14").Value = "India"

# Clear the clipboard marching-ants selection artifact.
$excel.CutCopyMode = $false

# Reflect the final cursor position in the sheet view.
$ws.Activate()
[void]$ws.Range("K14").Select()
